$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.023.78'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '2.949.93'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '379.41'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.63%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '101.36'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("E7").Value = '  +0.22%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '36.22'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("E11").Value = '  -0.48%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0848'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").Value = '3.412.28'
$ws.Range("E13").Value = '  -0.40%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '18.34'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("E15").Value = '  +4.42%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '11.89'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +66.79%  '
$ws.Range("D17").Value = '2.944.34'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("D19").Value = '51.003.32'
$ws.Range("E19").Value = '  -0.30%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '3.08'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.69%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.44'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").Value = '0.0₃0951'
$ws.Range("E22").Value = '  -0.27%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '69.47'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.45%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '266.76'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.79%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '3.21'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +11.02%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.13'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -3.30%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -8.72%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '25.62'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.38%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.163'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.84%  '
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("E32").Value = '  +2.88%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '50.53'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.15%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '33.54'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.19%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0431'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -5.28%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  +3.61%  '
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '16.70'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.50%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("E42").Value = '  +1.39%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '117.77'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.10%  '
$ws.Range("E44").Value = '  +8.73%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '21.36'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("E47").Value = '  -2.08%  '
$ws.Range("D48").Value = '2.008.69'
$ws.Range("E48").Value = '  -0.06%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.261'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -4.45%  '
$ws.Range("E50").Value = '  -9.65%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '5.29'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +4.21%  '
